$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C21: 1 -> 0.5
$ws.Range("C21").Value = 0.5

# Insert a new row at 22, pushing the former row 22 (Company 3 / Asset11) down to row 23
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the agriculture / AGRI11 asset data
$ws.Range("A22").Value = "AGRItest"
$ws.Range("B22").Value = "AGRI11"
$ws.Range("C22").Value = 0.5
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = "Agriculture"
$ws.Range("H22").Value = "Amazonas"
$ws.Range("L22").Value = "agriculture"
$ws.Range("M22").Value = "Soybean"

# Match the row height used for this new row in the target workbook
$ws.Rows.Item(22).RowHeight = 17.25
